$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '335.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.04%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.83%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.786'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.65%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08329'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.77%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '8.850'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.14%'

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.530'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.99%'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.976'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.47%'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.889'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.90%'

# Row 10
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9472'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.77%'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1240'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.03%'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1980'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.33%'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09839'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '6.06%'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04583'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '17.70%'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1068'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.82%'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001296'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.14%'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006066'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.37%'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.498'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.48%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.728'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.88%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1372'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.01%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2694'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '11.70%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04415'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.27%'

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.35%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004356'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.12%'

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '5.05%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003997'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-94.67%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02799'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '0.35%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05778'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '7.05%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007952'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.37%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1428'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.87%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008988'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.48%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002174'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.03%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01011'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-12.08%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007294'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '10.77%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.06%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003197'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.08%'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002274'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.34%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.06%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.06%'
